# Auto-generated edit script: updates cryptos list (rows 2-51, columns B-E)
# Commit: "Updated cryptos list on Thu May 25 06:42:48 UTC 2023 with GitHub Actions"
#
# All data cells in this sheet are stored as TEXT (not numbers), even price
# columns like "307.60" or "1.010" where Excel would otherwise auto-detect a
# number and silently drop the significant trailing zero. We therefore force
# the NumberFormat to Text ("@") before writing each value, then reset the
# cell style back to Normal afterwards so no stray style/format id is left
# behind (matching the original, unstyled data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "26.416.42"
Set-TextValue "E2" "  -1.63%  "
# Row 3
Set-TextValue "D3" "1.796.29"
Set-TextValue "E3" "  -1.84%  "
# Row 4
Set-TextValue "D4" "1.008"
Set-TextValue "E4" "  +0.13%  "
# Row 5
Set-TextValue "D5" "1.007"
Set-TextValue "E5" "  +0.17%  "
# Row 6
Set-TextValue "D6" "307.60"
Set-TextValue "E6" "  -0.82%  "
# Row 7
Set-TextValue "D7" "0.4541"
Set-TextValue "E7" "  -1.62%  "
# Row 8
Set-TextValue "D8" "0.3589"
Set-TextValue "E8" "  -2.18%  "
# Row 9
Set-TextValue "B9" "OKB"
Set-TextValue "C9" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D9" "46.26"
Set-TextValue "E9" "  +1.00%  "
# Row 10
Set-TextValue "B10" "Dogecoin"
Set-TextValue "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D10" "0.07123"
Set-TextValue "E10" "  -0.52%  "
# Row 11
Set-TextValue "B11" "Polygon"
Set-TextValue "C11" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D11" "0.8894"
Set-TextValue "E11" "  +1.58%  "
# Row 12
Set-TextValue "B12" "TRON"
Set-TextValue "C12" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.07815"
Set-TextValue "E12" "  -1.02%  "
# Row 13
Set-TextValue "B13" "Solana"
Set-TextValue "C13" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "19.52"
Set-TextValue "E13" "  -0.22%  "
# Row 14
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.775.26"
Set-TextValue "E14" "  -5.08%  "
# Row 15
Set-TextValue "B15" "Polkadot"
Set-TextValue "C15" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D15" "5.285"
Set-TextValue "E15" "  -0.89%  "
# Row 16
Set-TextValue "B16" "Chainlink"
Set-TextValue "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "6.336"
Set-TextValue "E16" "  -0.69%  "
# Row 17
Set-TextValue "B17" "Litecoin"
Set-TextValue "C17" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D17" "85.09"
Set-TextValue "E17" "  -3.01%  "
# Row 18
Set-TextValue "B18" "BinanceUSD"
Set-TextValue "C18" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D18" "1.010"
Set-TextValue "E18" "  +0.27%  "
# Row 19
Set-TextValue "B19" "ShibaInu"
Set-TextValue "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.000008591"
Set-TextValue "E19" "  -1.46%  "
# Row 20
Set-TextValue "B20" "Dai"
Set-TextValue "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.007"
Set-TextValue "E20" "  +0.15%  "
# Row 21
Set-TextValue "B21" "Avalanche"
Set-TextValue "C21" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D21" "14.29"
Set-TextValue "E21" "  -0.96%  "
# Row 22
Set-TextValue "B22" "WrappedBTC"
Set-TextValue "C22" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D22" "26.419.31"
Set-TextValue "E22" "  -1.72%  "
# Row 23
Set-TextValue "B23" "Uniswap"
Set-TextValue "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D23" "4.995"
Set-TextValue "E23" "  -0.06%  "
# Row 24
Set-TextValue "E24" "  +1.19%  "
# Row 25
Set-TextValue "B25" "WrappedliquidstakedEther2.0"
Set-TextValue "C25" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D25" "1.979.73"
Set-TextValue "E25" "  -5.41%  "
# Row 26
Set-TextValue "B26" "Toncoin"
Set-TextValue "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "1.978"
Set-TextValue "E26" "  -0.46%  "
# Row 27
Set-TextValue "B27" "Monero"
Set-TextValue "C27" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D27" "152.70"
Set-TextValue "E27" "  +1.31%  "
# Row 28
Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "17.92"
Set-TextValue "E28" "  -1.55%  "
# Row 29
Set-TextValue "B29" "LidoDAOToken"
Set-TextValue "C29" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D29" "2.047"
Set-TextValue "E29" "  +4.26%  "
# Row 30
Set-TextValue "B30" "BitcoinCash"
Set-TextValue "C30" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D30" "112.09"
Set-TextValue "E30" "  -1.21%  "
# Row 31
Set-TextValue "B31" "InternetComputer(DFINITY)"
Set-TextValue "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "4.872"
Set-TextValue "E31" "  -1.24%  "
# Row 32
Set-TextValue "B32" "Stellar"
Set-TextValue "C32" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D32" "0.08659"
Set-TextValue "E32" "  -2.04%  "
# Row 33
Set-TextValue "B33" "HuobiToken"
Set-TextValue "C33" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D33" "3.058"
Set-TextValue "E33" "  -2.13%  "
# Row 34
Set-TextValue "B34" "Filecoin"
Set-TextValue "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "4.455"
Set-TextValue "E34" "  +0.08%  "
# Row 35
Set-TextValue "B35" "ImmutableX"
Set-TextValue "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D35" "0.7290"
Set-TextValue "E35" "  -3.28%  "
# Row 36
Set-TextValue "B36" "RenderToken"
Set-TextValue "C36" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D36" "2.727"
Set-TextValue "E36" "  +6.22%  "
# Row 37
Set-TextValue "B37" "ARBITRUM"
Set-TextValue "C37" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D37" "1.112"
Set-TextValue "E37" "  -1.22%  "
# Row 38
Set-TextValue "B38" "TrustWalletToken"
Set-TextValue "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.076"
Set-TextValue "E38" "  -0.84%  "
# Row 39
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01942"
Set-TextValue "E39" "  +0.39%  "
# Row 40
Set-TextValue "B40" "Hedera"
Set-TextValue "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.05115"
Set-TextValue "E40" "  -0.19%  "
# Row 41
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.871"
Set-TextValue "E41" "  -1.81%  "
# Row 42
Set-TextValue "B42" "TheSandbox"
Set-TextValue "C42" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D42" "0.5139"
Set-TextValue "E42" "  +3.53%  "
# Row 43
Set-TextValue "B43" "FraxShare"
Set-TextValue "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "6.888"
Set-TextValue "E43" "  -0.06%  "
# Row 44
Set-TextValue "B44" "Algorand"
Set-TextValue "C44" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D44" "0.1519"
Set-TextValue "E44" "  -4.73%  "
# Row 45
Set-TextValue "B45" "Aptos"
Set-TextValue "C45" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D45" "8.012"
Set-TextValue "E45" "  -3.46%  "
# Row 46
Set-TextValue "D46" "1.008"
Set-TextValue "E46" "  +0.30%  "
# Row 47
Set-TextValue "B47" "Decentraland"
Set-TextValue "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D47" "0.4666"
Set-TextValue "E47" "  -0.11%  "
# Row 48
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.933"
Set-TextValue "E48" "  -1.41%  "
# Row 49
Set-TextValue "B49" "Quant"
Set-TextValue "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "100.47"
Set-TextValue "E49" "  -1.71%  "
# Row 50
Set-TextValue "B50" "NEARProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D50" "1.587"
Set-TextValue "E50" "  -1.39%  "
# Row 51
Set-TextValue "B51" "Cronos"
Set-TextValue "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.05979"
Set-TextValue "E51" "  -1.90%  "

Write-Output "Updated 181 cells in cryptos sheet"
